$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the header / identity fields of the bill form
$ws.Range("A3").Value2 = "নাম: Dr. M. M. A. Hashem"
$ws.Range("A4").Value2 = "পদবী: অধ্যাপক"
$ws.Range("F5").Value2 = "বিভাগ :সিএসই"

# Fill in the quantities for each billed item (label wise bill generate)
$ws.Range("G9").Value2 = 117
$ws.Range("G12").Value2 = 117
$ws.Range("G14").Value2 = 119
$ws.Range("G16").Value2 = 27
$ws.Range("G20").Value2 = 40
$ws.Range("G26").Value2 = 1
$ws.Range("G29").Value2 = 15

# Fill in the grand total amount written in words
$ws.Range("A32").Value2 = "কথায়:চার লক্ষ বত্রিশ হাজার একশত সাতান্ন টাকা মাত্র।"

# Leave the selection on the name field as in the final saved state
$ws.Range("B5").Select()
